# Daily attendance processing - 2025-10-06 20:20:13
# For every row in column G ("Recorded By"), rotate a comma-separated list of
# recorders so that the last entry moves to the front.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
